$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target cluster-pair table (Sending cluster, Ligand symbol, Receptor symbol, Target cluster, ...)
# now covers every combination of sending cluster {FAPs, sCs} x target cluster {ECs, FAPs, sCs}.

$data = @(
    @("FAPs","Lgi1","Adam23","ECs",2,0.6666666666666666,0.3953366666666667,1.18601,0.5534398610532347,0.5534398610532348,3,1,0.1145763333333333,0.343729,0.006557053879060051,0.006557053879060051,0.04529622569888889,0.40766603129,0.003628934987745568,0.003628934987745569),
    @("FAPs","Lgi1","Adam23","FAPs",2,0.6666666666666666,0.3953366666666667,1.18601,0.5534398610532347,0.5534398610532348,3,1,9.390663666666667,28.171991,0.5374154140831726,0.5374154140831726,3.712473671767778,33.41226304591,0.2974271120980577,0.2974271120980577),
    @("FAPs","Lgi1","Adam23","sCs",2,0.6666666666666666,0.3953366666666667,1.18601,0.5534398610532347,0.5534398610532348,3,1,7.968512,23.905536,0.4560275320377672,0.4560275320377672,3.150244972373333,28.35220475136,0.2523838139674315,0.2523838139674315),
    @("sCs","Lgi1","Adam23","ECs",3,1,0.3189896666666667,0.956969,0.4465601389467652,0.4465601389467652,3,1,0.1145763333333333,0.343729,0.006557053879060051,0.006557053879060051,0.03654866637788889,0.3289379974009999,0.002928118891314482,0.002928118891314482),
    @("sCs","Lgi1","Adam23","FAPs",3,1,0.3189896666666667,0.956969,0.4465601389467652,0.4465601389467652,3,1,9.390663666666667,28.171991,0.5374154140831726,0.5374154140831726,2.995524672808778,26.959722055279,0.2399883019851149,0.2399883019851149),
    @("sCs","Lgi1","Adam23","sCs",3,1,0.3189896666666667,0.956969,0.4465601389467652,0.4465601389467652,3,1,7.968512,23.905536,0.4560275320377672,0.4560275320377672,2.541872986709333,22.876856880384,0.2036437180703358,0.2036437180703358)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $col = $c + 1
        $ws.Cells.Item($row, $col).Value = $rowData[$c]
    }
}
